$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '70.788.50'
Set-TextCell $ws.Range("E2") '  +1.64%  '
Set-TextCell $ws.Range("D3") '3.632.80'
Set-TextCell $ws.Range("E3") '  +3.63%  '
Set-TextCell $ws.Range("E4") '  +0.06%  '
Set-TextCell $ws.Range("D5") '606.40'
Set-TextCell $ws.Range("E5") '  +0.22%  '
Set-TextCell $ws.Range("D6") '199.44'
Set-TextCell $ws.Range("E6") '  +2.38%  '
Set-TextCell $ws.Range("D7") '0.627'
Set-TextCell $ws.Range("E7") '  +0.39%  '
Set-TextCell $ws.Range("E8") '  +0.10%  '
Set-TextCell $ws.Range("E9") '  +10.27%  '
Set-TextCell $ws.Range("D10") '0.648'
Set-TextCell $ws.Range("E10") '  -0.14%  '
Set-TextCell $ws.Range("D11") '53.99'
Set-TextCell $ws.Range("E11") '  +1.14%  '
Set-TextCell $ws.Range("E12") '  +2.01%  '
Set-TextCell $ws.Range("E13") '  +0.78%  '
Set-TextCell $ws.Range("D14") '4.209.83'
Set-TextCell $ws.Range("E14") '  +3.60%  '
Set-TextCell $ws.Range("D15") '680.98'
Set-TextCell $ws.Range("E15") '  +14.76%  '
Set-TextCell $ws.Range("E16") '  +2.44%  '
Set-TextCell $ws.Range("D17") '70.902.29'
Set-TextCell $ws.Range("E17") '  +1.57%  '
Set-TextCell $ws.Range("D18") '3.655.11'
Set-TextCell $ws.Range("E18") '  +4.74%  '
Set-TextCell $ws.Range("D19") '19.02'
Set-TextCell $ws.Range("E19") '  -0.24%  '
Set-TextCell $ws.Range("E20") '  +0.38%  '
Set-TextCell $ws.Range("E21") '  +1.35%  '
Set-TextCell $ws.Range("D22") '18.67'
Set-TextCell $ws.Range("E22") '  +2.50%  '
Set-TextCell $ws.Range("D23") '5.38'
Set-TextCell $ws.Range("E23") '  +1.89%  '
Set-TextCell $ws.Range("D24") '105.95'
Set-TextCell $ws.Range("E24") '  +4.46%  '
Set-TextCell $ws.Range("E25") '  -0.27%  '
Set-TextCell $ws.Range("E26") '  -4.62%  '
Set-TextCell $ws.Range("D27") '10.45'
Set-TextCell $ws.Range("E27") '  -3.53%  '
Set-TextCell $ws.Range("D28") '9.85'
Set-TextCell $ws.Range("E28") '  +3.66%  '
Set-TextCell $ws.Range("D29") '34.24'
Set-TextCell $ws.Range("E29") '  +3.27%  '
Set-TextCell $ws.Range("D30") '4.66'
Set-TextCell $ws.Range("E30") '  +8.84%  '
Set-TextCell $ws.Range("E31") '  +1.66%  '
Set-TextCell $ws.Range("D32") '12.21'
Set-TextCell $ws.Range("E32") '  -1.39%  '
Set-TextCell $ws.Range("D33") '0.115'
Set-TextCell $ws.Range("E33") '  +0.56%  '
Set-TextCell $ws.Range("D34") '63.30'
Set-TextCell $ws.Range("E34") '  +0.31%  '
Set-TextCell $ws.Range("D35") '3.959.14'
Set-TextCell $ws.Range("E35") '  +6.30%  '
Set-TextCell $ws.Range("D36") '0.0₃0867'
Set-TextCell $ws.Range("E36") '  +5.59%  '
Set-TextCell $ws.Range("E37") '  -0.05%  '
Set-TextCell $ws.Range("D39") '36.83'
Set-TextCell $ws.Range("E39") '  +1.51%  '
Set-TextCell $ws.Range("D40") '505.08'
Set-TextCell $ws.Range("E40") '  +4.03%  '
Set-TextCell $ws.Range("E41") '  -0.44%  '
Set-TextCell $ws.Range("E42") '  -2.68%  '
Set-TextCell $ws.Range("E43") '  +2.66%  '
Set-TextCell $ws.Range("D44") '3.12'
Set-TextCell $ws.Range("E44") '  +11.03%  '
Set-TextCell $ws.Range("E45") '  +1.88%  '
Set-TextCell $ws.Range("E46") '  +6.53%  '
Set-TextCell $ws.Range("E47") '  +0.65%  '
Set-TextCell $ws.Range("E48") '  +3.52%  '
Set-TextCell $ws.Range("E49") '  -0.30%  '
Set-TextCell $ws.Range("E50") '  +0.90%  '
Set-TextCell $ws.Range("E51") '  +1.67%  '
